# Commit: "Added excel code and changed around writing"
#
# - Column E header (E1) held a stray "Unnamed: 4" label from a pandas
#   export; clear it out so the cell is blank instead.
# - The rest of column E (and, further down the sheet, columns D/E and
#   C/D/E) only ever held empty placeholder cells - clear those ranges so
#   the cells are dropped entirely rather than persisted as empty strings.
# - Give the sheet sane column widths instead of the default, and leave
#   the selection on E1 with the tab marked active, matching how the file
#   was left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blank out the leftover pandas "Unnamed: 4" header label.
$ws.Range("E1").Value = $null

# Drop the empty placeholder cells that trail the real data in each
# section of the sheet.
$ws.Range("E2:E26").ClearContents()
$ws.Range("D27:E52").ClearContents()
$ws.Range("C53:E78").ClearContents()

# Set explicit column widths (stored widths of 60 / 30 / 30 / 30 chars;
# Excel's ColumnWidth property is offset by ~5/6 of a character from the
# width actually persisted to the sheet, so back that out here).
$ws.Columns.Item(1).ColumnWidth = 59.1666666666667
$ws.Columns.Item(2).ColumnWidth = 29.1666666666667
$ws.Columns.Item(3).ColumnWidth = 29.1666666666667
$ws.Columns.Item(4).ColumnWidth = 29.1666666666667

# Leave the selection on E1 with the sheet tab active. (Wrapped in [void]
# so PowerShell doesn't echo Select()'s boolean return value.)
[void]$ws.Range("E1").Select()
